$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vtab = [char]11

$t.Cell(1,1).Range.Text = "55 x 85" + $vtab + "  8    5" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "5|    |"
$t.Cell(1,2).Range.Text = "32 x 90" + $vtab + "  9    0" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "2|    |"
$t.Cell(1,3).Range.Text = "47 x 31" + $vtab + "  3    1" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "7|    |"
$t.Cell(2,1).Range.Text = "93 x 13" + $vtab + "  1    3" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "3|    |"
$t.Cell(2,2).Range.Text = "33 x 99" + $vtab + "  9    9" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "3|    |"
$t.Cell(2,3).Range.Text = "49 x 14" + $vtab + "  1    4" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "9|    |"
$t.Cell(3,1).Range.Text = "43 x 43" + $vtab + "  4    3" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "3|    |"
$t.Cell(3,2).Range.Text = "21 x 43" + $vtab + "  4    3" + $vtab + "  ----" + $vtab + "2|    |" + $vtab + "1|    |"
$t.Cell(3,3).Range.Text = "17 x 91" + $vtab + "  9    1" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "7|    |"
$t.Cell(4,1).Range.Text = "15 x 95" + $vtab + "  9    5" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "5|    |"
$t.Cell(4,2).Range.Text = "98 x 54" + $vtab + "  5    4" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "8|    |"
$t.Cell(4,3).Range.Text = "24 x 20" + $vtab + "  2    0" + $vtab + "  ----" + $vtab + "2|    |" + $vtab + "4|    |"
$t.Cell(5,1).Range.Text = "16 x 64" + $vtab + "  6    4" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "6|    |"
$t.Cell(5,2).Range.Text = "81 x 41" + $vtab + "  4    1" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "1|    |"
$t.Cell(5,3).Range.Text = "69 x 74" + $vtab + "  7    4" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "9|    |"
